$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.832.47"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.637.67"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'216.87"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "'0.0622"
$ws.Range("D10").Value = "'19.89"
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "1.868.50"
$ws.Range("D13").Value = "1.640.41"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "'4.11"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'0.529"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "'66.76"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "26.832.42"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "'219.10"
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'6.73"
$ws.Range("E21").Value = "  +6.49%  "
$ws.Range("D22").Value = "'4.39"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +3.50%  "
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'147.01"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "'7.35"
$ws.Range("E27").Value = "  +4.72%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'15.78"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").Value = "'0.0502"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'3.33"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").Value = "1.262.67"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").Value = "'0.532"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "'0.831"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D41").Value = "'0.806"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'5.43"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "1.783.08"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "'61.66"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").Value = "'92.03"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  +17.37%  "
$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "'7.62"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'0.0960"
$ws.Range("E51").Value = "  -0.08%  "
